# Auto-generated edit script: update cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.397.71"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "3.414.18"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'570.35"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'156.68"
$ws.Range("E6").Value = "  -2.77%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +8.09%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "3.417.02"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").Value = "'7.13"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").Value = "4.003.71"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("E15").Value = "  -3.18%  "

$ws.Range("D16").Value = "'27.68"
$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("D17").Value = "64.441.28"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "3.416.52"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").Value = "'6.30"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").Value = "'13.86"
$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("D21").Value = "'377.71"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D23").Value = "'0.547"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").Value = "'71.75"
$ws.Range("E25").Value = "  -2.06%  "

$ws.Range("E26").Value = "  -3.99%  "

$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +4.51%  "

$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("E32").Value = "  -1.86%  "

$ws.Range("D33").Value = "'23.05"
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +6.92%  "

$ws.Range("D36").Value = "'160.02"
$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +5.67%  "

$ws.Range("D39").Value = "'0.0763"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").Value = "2.875.57"
$ws.Range("E40").Value = "  -4.31%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.64"
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'26.42"
$ws.Range("E42").Value = "  -3.17%  "

$ws.Range("D43").Value = "'42.84"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'26.32"
$ws.Range("E44").Value = "  +6.08%  "

$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "'0.769"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").Value = "'321.35"
$ws.Range("E47").Value = "  +5.75%  "

$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("E49").Value = "  +2.49%  "

$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("E51").Value = "  -2.37%  "
